$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$fValues = @{
    2  = -3
    4  = 4
    5  = 6
    6  = -6
    7  = 1
    9  = 3
    10 = -4
    11 = -3
    12 = -2
    13 = -4
    14 = 1
    15 = 1
    16 = 3
    17 = 7
    19 = 1
    20 = 2
    21 = -6
    22 = 4
    23 = -2
    24 = 1
    25 = 0
    26 = 1
    27 = 5
    28 = -1
    30 = -2
    31 = -2
    32 = 3
    34 = 1
    35 = -1
    36 = -2
    38 = -2
    39 = 1
}

foreach ($row in $fValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $fValues[$row]
}
